$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Agtr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.788972999999999
$ws.Range("H2").Value = 17.366919
$ws.Range("I2").Value = 0.9474260381515466
$ws.Range("J2").Value = 0.9474260381515467
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.202916333333333
$ws.Range("N2").Value = 3.608749
$ws.Range("O2").Value = 0.02173550655091635
$ws.Range("P2").Value = 0.02173550655091635
$ws.Range("Q2").Value = 6.963650174925667
$ws.Range("R2").Value = 62.67285157433101
$ws.Range("S2").Value = 0.02059278485875167
$ws.Range("T2").Value = 0.02059278485875167

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Agtr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.788972999999999
$ws.Range("H3").Value = 17.366919
$ws.Range("I3").Value = 0.9474260381515466
$ws.Range("J3").Value = 0.9474260381515467
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.52978033333333
$ws.Range("N3").Value = 61.589341
$ws.Range("O3").Value = 0.3709527941045833
$ws.Range("P3").Value = 0.3709527941045833
$ws.Range("Q3").Value = 118.8463440455977
$ws.Range("R3").Value = 1069.617096410379
$ws.Range("S3").Value = 0.3514503360597517
$ws.Range("T3").Value = 0.3514503360597518

$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Agtr1a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.788972999999999
$ws.Range("H4").Value = 17.366919
$ws.Range("I4").Value = 0.9474260381515466
$ws.Range("J4").Value = 0.9474260381515467
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 33.61068033333333
$ws.Range("N4").Value = 100.832041
$ws.Range("O4").Value = 0.6073116993445002
$ws.Range("P4").Value = 0.6073116993445002
$ws.Range("Q4").Value = 194.5713209612977
$ws.Range("R4").Value = 1751.141888651679
$ws.Range("S4").Value = 0.575382917233043
$ws.Range("T4").Value = 0.5753829172330431

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Agt"
$ws.Range("C5").Value = "Agtr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.321238
$ws.Range("H5").Value = 0.963714
$ws.Range("I5").Value = 0.05257396184845335
$ws.Range("J5").Value = 0.05257396184845335
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.202916333333333
$ws.Range("N5").Value = 3.608749
$ws.Range("O5").Value = 0.02173550655091635
$ws.Range("P5").Value = 0.02173550655091635
$ws.Range("Q5").Value = 0.3864224370873333
$ws.Range("R5").Value = 3.477801933786
$ws.Range("S5").Value = 0.001142721692164684
$ws.Range("T5").Value = 0.001142721692164684

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Agt"
$ws.Range("C6").Value = "Agtr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.321238
$ws.Range("H6").Value = 0.963714
$ws.Range("I6").Value = 0.05257396184845335
$ws.Range("J6").Value = 0.05257396184845335
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 20.52978033333333
$ws.Range("N6").Value = 61.589341
$ws.Range("O6").Value = 0.3709527941045833
$ws.Range("P6").Value = 0.3709527941045833
$ws.Range("Q6").Value = 6.594945574719334
$ws.Range("R6").Value = 59.354510172474
$ws.Range("S6").Value = 0.01950245804483153
$ws.Range("T6").Value = 0.01950245804483153

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Agt"
$ws.Range("C7").Value = "Agtr1a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.321238
$ws.Range("H7").Value = 0.963714
$ws.Range("I7").Value = 0.05257396184845335
$ws.Range("J7").Value = 0.05257396184845335
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 33.61068033333333
$ws.Range("N7").Value = 100.832041
$ws.Range("O7").Value = 0.6073116993445002
$ws.Range("P7").Value = 0.6073116993445002
$ws.Range("Q7").Value = 10.79702772891933
$ws.Range("R7").Value = 97.173249560274
$ws.Range("S7").Value = 0.03192878211145713
$ws.Range("T7").Value = 0.03192878211145713

